$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1.284767783027405, 0.2638960490923239, 0.2453400986128003, 1.667739334718291, 0.9777017153308094, 0.9546698591763203, 0.3492363118898965, 1.115625055482166),
    @(3, 1.166196123525538, 0.2323419556084332, 0.2400956524805054, 1.647182596362242, 0.9600047766263344, 0.9527192173408849, 0.3379428821024959, 1.132091594119473),
    @(4, 1.093786557672786, 0.2129811081162245, 0.2369705582724606, 1.635756539032968, 0.949973871096347, 0.9521388886133764, 0.3312283164755172, 1.142736083101809),
    @(5, 1.064378089576508, 0.2050947716413134, 0.2357209969206622, 1.631399869335127, 0.9460949653806239, 0.952057214960206, 0.328547053159582, 1.147207913221489),
    @(6, 1.059500815836202, 0.2037854541345325, 0.2355149555389602, 1.630694505255192, 0.9454634496840413, 0.9520529911085163, 0.3281051444963339, 1.147958552934494),
    @(7, 1.093389543902276, 0.2128747366235473, 0.2369536092544138, 1.635696572182468, 0.9499207151335014, 0.9521371608688725, 0.3311919337547948, 1.142795848954307),
    @(8, 1.243802453991123, 0.253013369832729, 0.2435121186358202, 1.660402453313665, 0.9714257447774912, 0.9538688796606607, 0.3452965912638462, 1.121191799347397),
    @(9, 1.541905582293566, 0.3318376098266356, 0.2571261240847065, 1.718397135132378, 1.020281844194045, 0.9621843780062136, 0.3747119994741297, 1.083069679409995),
    @(10, 1.762888176602587, 0.3898319845499145, 0.2675871360690394, 1.766913295415279, 1.060338430799675, 0.9713244547019428, 0.3974155006022357, 1.0576566177925),
    @(11, 1.863858268641138, 0.4162363507829809, 0.2724458327740535, 1.790286827347458, 1.079484097221268, 0.9761476343954598, 0.4079858486428236, 1.046661003886914),
    @(12, 1.902157381051666, 0.4262383725887844, 0.2743000448542432, 1.799326646775, 1.086868354273776, 0.978070256937059, 0.4120237623592686, 1.042578687994244),
    @(13, 1.893906137088493, 0.4240841125620136, 0.2739000704495567, 1.797371342519966, 1.085272031364013, 0.9776518994225682, 0.4111525572632928, 1.043454259649842),
    @(14, 1.8670078778689, 0.4170591575513072, 0.2725980929707248, 1.791026747662983, 1.080088906818077, 0.9763038789491247, 0.408317344032298, 1.046323514770865),
    @(15, 1.850540245804893, 0.4127565956366652, 0.271802459056218, 1.787165124137289, 1.076931612397459, 0.9754907196145268, 0.4065852807223962, 1.048091635622058),
    @(16, 1.756298484695435, 0.3881068514492085, 0.2672716161445408, 1.765412121459647, 1.059105921738876, 0.9710226804930358, 0.3967296064661809, 1.058386595781975),
    @(17, 1.698597988685492, 0.3729908011588918, 0.2645176557878415, 1.752402124568135, 1.048408003096256, 0.9684524523579228, 0.3907457647569572, 1.064847084387807),
    @(18, 1.665451974764892, 0.3642985579976994, 0.2629430588493165, 1.745041707949852, 1.042341705508477, 0.9670367162020455, 0.3873268163838191, 1.068616143625693),
    @(19, 1.654236478449263, 0.361355875085053, 0.2624115445158139, 1.742570611021222, 1.040302638227928, 0.9665681066560126, 0.3861731260917907, 1.069901404589437),
    @(20, 1.704735979835789, 0.3745997119794424, 0.2648098461533692, 1.75377436363155, 1.049537814952515, 0.9687195758825169, 0.3913803931809383, 1.064153851026404),
    @(21, 1.874906812821052, 0.4191224683044652, 0.2729801267454803, 1.792885174898132, 1.081607662714248, 0.976697210536571, 0.4091491584992042, 1.045478531491334),
    @(22, 1.986495975228195, 0.4482396934269559, 0.2784033804464343, 1.819547288754976, 1.103350080518851, 0.9824719641947013, 0.4209671076840351, 1.033748206014725),
    @(23, 1.926904605493803, 0.4326975280156375, 0.2755012604182241, 1.805216028179686, 1.091673647612822, 0.9793383736185319, 0.4146407870390618, 1.039965347328579),
    @(24, 1.701960910660375, 0.3738723290874759, 0.2646777198141876, 1.753153603144639, 1.049026765222294, 0.9685986164905387, 0.3910934115307896, 1.064467091027474),
    @(25, 1.460919158206309, 0.3105001317534857, 0.2533626010243779, 1.701677004458588, 1.0063403881515, 0.9594048973478948, 0.3665641422305583, 1.092927603937252)
)

$cols = @("B", "C", "D", "F", "G", "H", "J", "N")

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $row[$i + 1]
        $addr = $col + $r
        $ws.Range($addr).Value = $val
    }
}

$wb.Save()